$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 525.5714
$ws.Range("I4").Value = 535.8
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 535.8
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -421.8
$ws.Range("N4").Value = -728
$ws.Range("H34").Value = 6094.615
$ws.Range("I34").Value = 2536.889
$ws.Range("J34").Value = 14099.5
$ws.Range("K34").Value = 2536.889
$ws.Range("L34").Value = 14099.5
$ws.Range("M34").Value = -2333.889
$ws.Range("N34").Value = -14505.5
$ws.Range("H36").Value = 6094.615
$ws.Range("I36").Value = 2536.889
$ws.Range("J36").Value = 14099.5
$ws.Range("K36").Value = 2536.889
$ws.Range("L36").Value = 14099.5
$ws.Range("M36").Value = -1821.889
$ws.Range("N36").Value = -15529.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 3250
$ws.Range("J10").Value = 2500
$ws.Range("L10").Value = 2500
$ws.Range("N10").Value = -2840
$ws.Range("I45").Value = 16281.714
$ws.Range("K45").Value = 16281.714
$ws.Range("M45").Value = -15904.714
$ws.Range("H61").Value = 227198.95
$ws.Range("I61").Value = 6163.9614
$ws.Range("J61").Value = 529667.9
$ws.Range("K61").Value = 6163.9614
$ws.Range("L61").Value = 529667.9
$ws.Range("M61").Value = -5951.9614
$ws.Range("N61").Value = -530091.9
$ws.Range("H106").Value = 43700
$ws.Range("J106").Value = 43700
$ws.Range("L106").Value = 43700
$ws.Range("N106").Value = -46224
$ws.Range("H136").Value = 227198.95
$ws.Range("I136").Value = 6163.9614
$ws.Range("J136").Value = 529667.9
$ws.Range("K136").Value = 18491.8842
$ws.Range("L136").Value = 1589003.7
$ws.Range("M136").Value = -15941.8842
$ws.Range("N136").Value = -1594103.7

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H107").Value = 1353.2222
$ws.Range("I107").Value = 1148.5
$ws.Range("J107").Value = 2991
$ws.Range("K107").Value = 1148.5
$ws.Range("L107").Value = 2991
$ws.Range("M107").Value = 771.5
$ws.Range("N107").Value = -6831

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 43429.285
$ws.Range("I4").Value = 90001
$ws.Range("J4").Value = 35667.332
$ws.Range("K4").Value = 90001
$ws.Range("L4").Value = 35667.332
$ws.Range("M4").Value = -89889
$ws.Range("N4").Value = -35891.332
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 214659.52
$ws.Range("I58").Value = 1672.2667
$ws.Range("J58").Value = 314497.28
$ws.Range("K58").Value = 1672.2667
$ws.Range("L58").Value = 314497.28
$ws.Range("M58").Value = -1469.2667
$ws.Range("N58").Value = -314903.28
$ws.Range("H99").Value = 5578.6924
$ws.Range("I99").Value = 6293
$ws.Range("K99").Value = 6293
$ws.Range("M99").Value = -4795
$ws.Range("H126").Value = 5578.6924
$ws.Range("I126").Value = 6293
$ws.Range("K126").Value = 18879
$ws.Range("M126").Value = -16409
$ws.Range("H132").Value = 2545.75
$ws.Range("I132").Value = 1778.2307
$ws.Range("J132").Value = 3971.1428
$ws.Range("K132").Value = 5334.6921
$ws.Range("L132").Value = 11913.4284
$ws.Range("M132").Value = -2804.6921
$ws.Range("N132").Value = -16973.4284
$ws.Range("H134").Value = 254396.83
$ws.Range("I134").Value = 5063.75
$ws.Range("J134").Value = 628396.4399999999
$ws.Range("K134").Value = 15191.25
$ws.Range("L134").Value = 1885189.32
$ws.Range("M134").Value = -12656.25
$ws.Range("N134").Value = -1890259.32
$ws.Range("H136").Value = 214659.52
$ws.Range("I136").Value = 1672.2667
$ws.Range("J136").Value = 314497.28
$ws.Range("K136").Value = 5016.800099999999
$ws.Range("L136").Value = 943491.8400000001
$ws.Range("M136").Value = -2466.800099999999
$ws.Range("N136").Value = -948591.8400000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 213.6
$ws.Range("J10").Value = 333.33334
$ws.Range("L10").Value = 1000.00002
$ws.Range("N10").Value = -1278.00002
$ws.Range("H17").Value = 792.5
$ws.Range("I17").Value = 792.5
$ws.Range("K17").Value = 2377.5
$ws.Range("M17").Value = -2208.5
$ws.Range("H56").Value = 5990
$ws.Range("I56").Value = 5990
$ws.Range("K56").Value = 5990
$ws.Range("M56").Value = -5460

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5422.222
$ws.Range("J2").Value = 5422.222
$ws.Range("L2").Value = 5422.222
$ws.Range("N2").Value = -5646.222
$ws.Range("H10").Value = 2101
$ws.Range("I10").Value = 2550.75
$ws.Range("J10").Value = 1501.3334
$ws.Range("K10").Value = 2550.75
$ws.Range("L10").Value = 1501.3334
$ws.Range("M10").Value = -2410.75
$ws.Range("N10").Value = -1781.3334
$ws.Range("H22").Value = 1916.8334
$ws.Range("I22").Value = 501
$ws.Range("J22").Value = 2000.1177
$ws.Range("K22").Value = 501
$ws.Range("L22").Value = 2000.1177
$ws.Range("M22").Value = -206
$ws.Range("N22").Value = -2590.1177
$ws.Range("H27").Value = 1916.8334
$ws.Range("I27").Value = 501
$ws.Range("J27").Value = 2000.1177
$ws.Range("K27").Value = 501
$ws.Range("L27").Value = 2000.1177
$ws.Range("M27").Value = -394
$ws.Range("N27").Value = -2214.1177
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 1357
$ws.Range("I31").Value = 650.1111
$ws.Range("J31").Value = 2947.5
$ws.Range("K31").Value = 650.1111
$ws.Range("L31").Value = 2947.5
$ws.Range("M31").Value = -402.1111
$ws.Range("N31").Value = -3443.5
$ws.Range("H32").Value = 4945
$ws.Range("J32").Value = 8199
$ws.Range("L32").Value = 8199
$ws.Range("N32").Value = -8833
$ws.Range("H40").Value = 47621570
$ws.Range("I40").Value = 52634104
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 52634104
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -52633968
$ws.Range("N40").Value = -2772
$ws.Range("H46").Value = 733.6667
$ws.Range("I46").Value = 733.6667
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 733.6667
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -545.6667
$ws.Range("N46").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17387.875
$ws.Range("J2").Value = 19743.285
$ws.Range("L2").Value = 19743.285
$ws.Range("N2").Value = -19967.285
$ws.Range("H11").Value = 1374500
$ws.Range("I11").Value = 749000
$ws.Range("K11").Value = 749000
$ws.Range("M11").Value = -748858
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2939
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -14696
$ws.Range("N84").ClearContents()
$ws.Range("H117").Value = 24877
$ws.Range("J117").Value = 24877
$ws.Range("L117").Value = 24877
$ws.Range("N117").Value = -34055
